# Updated Kevin hour log
# Adds a new work-log entry (row 17) to Sheet1, describing work completed on
# drawing/wiring the IWRL6432AOP along with a clock oscillator and power
# sources, and updates the sheet view/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting (fonts/fills/borders/number formats) of the last
# existing data row (row 16) down into the new row 17 before writing values,
# so the new row matches the look of the rest of the log.
$ws.Range("B16:G16").Copy()
$ws.Range("B17:G17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row of data.
$ws.Range("B17").Value = "Kevin Su"
$ws.Range("C17").Value = 45682
$ws.Range("D17").Value = 0.45833333333333331
$ws.Range("E17").Value = 0.70833333333333337
$ws.Range("F17").Value = 0.25
$ws.Range("G17").Value = "Completed drawing IWRL6432AOP and wired it in the schematic along with an clock oscillator and power sources."

# Match the taller row height used for this entry's wrapped description text.
$ws.Rows.Item(17).RowHeight = 64

# Widen column G so the new (and existing) descriptions are easier to read.
$ws.Columns.Item(7).ColumnWidth = 24.83

# Scroll the view down and move the active selection the same way Excel
# would after typing the new row.
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1
[void]$ws.Range("H17:H23").Select()
